# finalize shipley OES lookup; continue evolution of generic method
#
# This script:
#  1. Adds a new "child" worksheet after the existing "adult" worksheet,
#     populated with the CV68 child lookup-table data.
#  2. Formats the two text-like age-strat codes ("084" / "096") as text
#     (right aligned) so the leading zero is preserved.
#  3. Sets reasonable column widths on the new sheet.
#  4. Updates the selection / active sheet so that "child" becomes the
#     active (visible) tab, matching the final state of the workbook.
#  5. Leaves the "adult" sheet's data untouched, only updating its
#     sheet view (no longer the selected tab, new selection rectangle).

$wb = $excel.ActiveWorkbook
$adult = $wb.Worksheets.Item("adult")

# ---------------------------------------------------------------------
# 1. Create the new "child" worksheet right after "adult"
# ---------------------------------------------------------------------
$child = $wb.Worksheets.Add($null, $adult)
$child.Name = "child"

# ---------------------------------------------------------------------
# 2. Header row (reuses the same shared-string header labels as "adult")
# ---------------------------------------------------------------------
$child.Range("A1").Value = "agestrat"
$child.Range("B1").Value = "VOC_CV68"
$child.Range("C1").Value = "ABS_CV68"
$child.Range("D1").Value = "BLO_CV68"
$child.Range("E1").Value = "CMA_CV68"
$child.Range("F1").Value = "CMB_CV68"

# ---------------------------------------------------------------------
# 3. Data rows
# ---------------------------------------------------------------------
# Rows whose "agestrat" is textual ("084", "096") - written with a
# leading apostrophe so Excel keeps them as text (preserving the
# leading zero), then normalized back to the "Normal" style and
# right-aligned, matching the final look of typed-then-reformatted
# text cells.
$child.Range("A2").Value = "'084"
$child.Range("B2").Value = 6
$child.Range("C2").Value = 8
$child.Range("D2").Value = 8
$child.Range("E2").Value = 6
$child.Range("F2").Value = 7

$child.Range("A3").Value = "'096"
$child.Range("B3").Value = 6
$child.Range("C3").Value = 7
$child.Range("D3").Value = 7
$child.Range("E3").Value = 5
$child.Range("F3").Value = 6

$child.Range("A2:A3").Style = "Normal"
$child.Range("A2:A3").HorizontalAlignment = -4152  # xlRight

# Remaining rows use plain numeric agestrat values
$child.Range("A4").Value = 108
$child.Range("B4").Value = 6
$child.Range("C4").Value = 8
$child.Range("D4").Value = 5
$child.Range("E4").Value = 6
$child.Range("F4").Value = 6

$child.Range("A5").Value = 120
$child.Range("B5").Value = 7
$child.Range("C5").Value = 7
$child.Range("D5").Value = 6
$child.Range("E5").Value = 6
$child.Range("F5").Value = 5

$child.Range("A6").Value = 132
$child.Range("B6").Value = 6
$child.Range("C6").Value = 8
$child.Range("D6").Value = 6
$child.Range("E6").Value = 5
$child.Range("F6").Value = 5

$child.Range("A7").Value = 144
$child.Range("B7").Value = 6
$child.Range("C7").Value = 8
$child.Range("D7").Value = 6
$child.Range("E7").Value = 6
$child.Range("F7").Value = 5

$child.Range("A8").Value = 156
$child.Range("B8").Value = 6
$child.Range("C8").Value = 7
$child.Range("D8").Value = 5
$child.Range("E8").Value = 5
$child.Range("F8").Value = 5

$child.Range("A9").Value = 180
$child.Range("B9").Value = 5
$child.Range("C9").Value = 7
$child.Range("D9").Value = 4
$child.Range("E9").Value = 4
$child.Range("F9").Value = 4

$child.Range("A10").Value = 204
$child.Range("B10").Value = 6
$child.Range("C10").Value = 7
$child.Range("D10").Value = 4
$child.Range("E10").Value = 5
$child.Range("F10").Value = 4

# ---------------------------------------------------------------------
# 4. Column widths for the new sheet
# ---------------------------------------------------------------------
$child.Columns.Item(2).ColumnWidth = 10.6640625
$child.Columns.Item(3).ColumnWidth = 9.6640625
$child.Range("D1:E1").EntireColumn.ColumnWidth = 10
$child.Columns.Item(6).ColumnWidth = 11.109375

# ---------------------------------------------------------------------
# 5. Selections / active tab
# ---------------------------------------------------------------------
# "adult" is no longer the selected tab; its selection becomes A1:F1
$adult.Range("A1:F1").Select()

# "child" becomes the active sheet/tab, with F11 selected (first empty
# cell below the data, as left by the person entering the table)
$child.Activate()
$child.Range("F11").Select()

Write-Host "Added 'child' worksheet with CV68 lookup data."
